# Update cryptocurrency price/volume snapshot cells (columns D and E).
# Values are kept as text (matching the source inline-string cells), so
# NumberFormat is forced to "@" before assignment and the style reset
# back to "Normal" afterwards to avoid leaving a lingering text format.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '68.398.68'),
    @('E2', '  +0.18%  '),
    @('D3', '2.647.23'),
    @('E3', '  +0.19%  '),
    @('D5', '596.71'),
    @('E5', '  -0.20%  '),
    @('D6', '158.84'),
    @('E6', '  +2.83%  '),
    @('E7', '  -0.01%  '),
    @('D8', '0.540'),
    @('E8', '  -1.02%  '),
    @('D9', '2.647.41'),
    @('E9', '  +0.18%  '),
    @('E10', '  -2.15%  '),
    @('E11', '  -1.01%  '),
    @('E12', '  +0.42%  '),
    @('D13', '0.351'),
    @('E13', '  -0.65%  '),
    @('E14', '  +0.26%  '),
    @('D15', '3.131.15'),
    @('E15', '  +0.25%  '),
    @('E16', '  -3.15%  '),
    @('D17', '68.297.89'),
    @('E17', '  +0.02%  '),
    @('D18', '2.643.97'),
    @('E18', '  -0.37%  '),
    @('E19', '  +1.70%  '),
    @('D20', '363.35'),
    @('E20', '  -0.29%  '),
    @('D21', '7.47'),
    @('E21', '  +0.24%  '),
    @('E22', '  +1.09%  '),
    @('E23', '  -1.88%  '),
    @('E24', '  +0.84%  '),
    @('D25', '74.74'),
    @('E25', '  -0.88%  '),
    @('E26', '  +0.10%  '),
    @('D27', '9.90'),
    @('E27', '  +1.28%  '),
    @('D28', '2.795.67'),
    @('E28', '  +0.60%  '),
    @('E29', '  -2.86%  '),
    @('E30', '  -0.17%  '),
    @('D31', '566.69'),
    @('E31', '  +1.00%  '),
    @('D32', '8.05'),
    @('E32', '  -0.07%  '),
    @('E33', '  -0.17%  '),
    @('E34', '  +0.50%  '),
    @('E35', '  +4.55%  '),
    @('E36', '  -1.55%  '),
    @('D37', '1.00'),
    @('E37', '  -0.01%  '),
    @('D38', '160.75'),
    @('E38', '  -0.47%  '),
    @('E39', '  +1.76%  '),
    @('E40', '  -1.15%  '),
    @('E41', '  -0.85%  '),
    @('D42', '5.31'),
    @('E42', '  -0.60%  '),
    @('E43', '  +0.47%  '),
    @('E44', '  -5.50%  '),
    @('E45', '  +0.12%  '),
    @('D46', '158.22'),
    @('E46', '  +1.26%  '),
    @('E47', '  +1.84%  '),
    @('D48', '21.92'),
    @('E48', '  +0.86%  '),
    @('E49', '  -0.09%  '),
    @('E50', '  -1.15%  '),
    @('D51', '0.574'),
    @('E51', '  +2.00%  ')
)

foreach ($pair in $updates) {
    $ref = $pair[0]
    $value = $pair[1]
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}
